# Apply quantity corrections to FullOrderLog and FullShipmentLog sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("FullOrderLog")
$ws2 = $wb.Worksheets.Item("FullShipmentLog")

# FullOrderLog sheet - column E ("quantity")
$ws1.Range("E27").Value = 8
$ws1.Range("E28").Value = 40
$ws1.Range("E30").Value = 100
$ws1.Range("E33").Value = 23
$ws1.Range("E40").Value = 8
$ws1.Range("E43").Value = 8
$ws1.Range("E44").Value = 42
$ws1.Range("E45").Value = 14
$ws1.Range("E46").Value = 100
$ws1.Range("E49").Value = 24
$ws1.Range("E53").Value = 8
$ws1.Range("E56").Value = 8
$ws1.Range("E57").Value = 14
$ws1.Range("E60").Value = 39
$ws1.Range("E61").Value = 15
$ws1.Range("E62").Value = 108
$ws1.Range("E65").Value = 21
$ws1.Range("E67").Value = 8
$ws1.Range("E69").Value = 8
$ws1.Range("E70").Value = 14
$ws1.Range("E73").Value = 15
$ws1.Range("E76").Value = 40
$ws1.Range("E77").Value = 14
$ws1.Range("E78").Value = 108
$ws1.Range("E81").Value = 23
$ws1.Range("E83").Value = 8
$ws1.Range("E84").Value = 14
$ws1.Range("E86").Value = 15
$ws1.Range("E89").Value = 14
$ws1.Range("E93").Value = 14
$ws1.Range("E94").Value = 106
$ws1.Range("E97").Value = 24

# FullShipmentLog sheet - column D ("quantity")
$ws2.Range("D14").Value = 25
$ws2.Range("D17").Value = 40
$ws2.Range("D19").Value = 4
$ws2.Range("D20").Value = 49
$ws2.Range("D21").Value = 114
$ws2.Range("D23").Value = 26
$ws2.Range("D26").Value = 46
$ws2.Range("D27").Value = 85
$ws2.Range("D29").Value = 23
$ws2.Range("D32").Value = 48
$ws2.Range("D33").Value = 104
$ws2.Range("D35").Value = 25
$ws2.Range("D38").Value = 49
$ws2.Range("D39").Value = 53
$ws2.Range("D41").Value = 27
